$d = $word.ActiveDocument

# The document has three picture logos living in the headers/footers:
#   - Section 1, Header(2) [first-page header]: BTec logo  -> rename image1.jpg to image2.jpg
#   - Section 1, Footer(1) [primary footer]:    Pearson logo -> rename image2.png to image1.png
#   - Section 1, Footer(2) [first-page footer]: Pearson logo -> rename image2.png to image1.png
#
# Each shape is re-fetched from ActiveDocument immediately before it is used so the
# handle can never go stale after an earlier header/footer edit in this same script.

$d.Sections.Item(1).Headers.Item(2).Range.InlineShapes.Item(1).Name = "image2.jpg"

$d.Sections.Item(1).Footers.Item(1).Range.InlineShapes.Item(1).Name = "image1.png"

$d.Sections.Item(1).Footers.Item(2).Range.InlineShapes.Item(1).Name = "image1.png"
